# Update cryptocurrency price/volume snapshot data (GitHub Actions scheduled refresh).
# For rows whose Price column previously held a plain decimal-looking string (e.g. "1.001"),
# the cell is pre-formatted as Text ("@") before the assignment so Excel keeps storing it
# as a string instead of silently converting it to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.710.41"
$ws.Range("E2").Value = "  +0.56%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "1.888.21"
$ws.Range("E3").Value = "  +0.28%  "
# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.55"
$ws.Range("E5").Value = "  +0.08%  "
# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").Value = "  -0.42%  "
# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2928"
$ws.Range("E8").Value = "  -0.17%  "
# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06532"
$ws.Range("E9").Value = "  +0.01%  "
# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  -0.21%  "
# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07794"
$ws.Range("E11").Value = "  +0.78%  "
# Row 12: Litecoin
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.892.82"
$ws.Range("E12").Value = "  +0.45%  "
# Row 13: WrappedEther
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.03"
$ws.Range("E13").Value = "  -1.13%  "
# Row 14: Polygon
$ws.Range("E14").Value = "  -0.43%  "
# Row 15: Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.250"
$ws.Range("E15").Value = "  +1.93%  "
# Row 16: BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.04"
$ws.Range("E16").Value = "  +3.58%  "
# Row 17: WrappedBTC
$ws.Range("D17").Value = "30.815.75"
$ws.Range("E17").Value = "  +0.92%  "
# Row 18: Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.18"
$ws.Range("E18").Value = "  -2.44%  "
# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007560"
$ws.Range("E19").Value = "  -0.20%  "
# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.141.28"
$ws.Range("E21").Value = "  +0.40%  "
# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.315"
$ws.Range("E22").Value = "  +1.07%  "
# Row 23: BinanceUSD
$ws.Range("E23").Value = "  +0.04%  "
# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.272"
# Row 25: Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.230"
$ws.Range("E25").Value = "  -0.72%  "
# Row 26: Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.06"
$ws.Range("E26").Value = "  +0.18%  "
# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("E27").Value = "  +0.08%  "
# Row 28: LidoDAOToken
$ws.Range("E28").Value = "  -1.21%  "
# Row 29: Toncoin
$ws.Range("E29").Value = "  -1.88%  "
# Row 30: Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09724"
$ws.Range("E30").Value = "  -3.42%  "
# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.495"
$ws.Range("E31").Value = "  -1.76%  "
# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.290"
$ws.Range("E32").Value = "  -1.01%  "
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.201"
$ws.Range("E33").Value = "  +2.15%  "
# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04851"
$ws.Range("E34").Value = "  +0.51%  "
# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.126"
$ws.Range("E35").Value = "  -0.51%  "
# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6983"
$ws.Range("E36").Value = "  -0.57%  "
# Row 37: HuobiToken
$ws.Range("E37").Value = "  +0.40%  "
# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01909"
$ws.Range("E38").Value = "  +2.17%  "
# Row 39: MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.811"
$ws.Range("E39").Value = "  +2.10%  "
# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.358"
$ws.Range("E40").Value = "  +0.61%  "
# Row 41: Aave
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.05"
$ws.Range("E41").Value = "  +6.05%  "
# Row 42: RenderToken
$ws.Range("E42").Value = "  +0.95%  "
# Row 43: TheSandbox
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4264"
$ws.Range("E43").Value = "  +0.95%  "
# Row 44: PaxDollar
$ws.Range("E44").Value = "  +0.02%  "
# Row 45: TrustWalletToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8365"
$ws.Range("E45").Value = "  -0.99%  "
# Row 46: Quant
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.37"
$ws.Range("E46").Value = "  -1.53%  "
# Row 47: EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.485"
$ws.Range("E47").Value = "  +1.91%  "
# Row 48: Aptos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.057"
$ws.Range("E48").Value = "  -0.69%  "
# Row 49: Elrond
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.66"
$ws.Range("E49").Value = "  +0.02%  "
# Row 50: Maker
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "917.50"
$ws.Range("E50").Value = "  -0.09%  "
# Row 51: Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05755"
$ws.Range("E51").Value = "  +1.89%  "
